# Automatische test-sync: 2025-06-19 19:50:30
#
# 1) Append two new "Afmelding nieuwsbrief" rows (43 & 44) to the Logs sheet.
# 2) Extend the conditional formatting ranges on the Logs sheet to cover the
#    new rows (D2:D42 -> D2:D44, G2:G42 -> G2:G44).
# 3) Re-sort/re-derive the Dashboard category summary (rows 6-13) to include
#    the new "Afmelding / Nieuwsbrief" category and keep it ordered by count.
# 4) Point the Dashboard bar chart's category/value series at the extended
#    Dashboard range (A2:A13 / B2:B13).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Logs sheet: add the two new rows
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(43, 1).Value = "Afmelding nieuwsbrief"
$logs.Cells.Item(43, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(43, 3).Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Cells.Item(43, 4).Value = "Afmelding / Nieuwsbrief"
$logs.Cells.Item(43, 6).Value = "2025-06-19 19:50:26"
$logs.Cells.Item(43, 7).Value = "Nee"

$logs.Cells.Item(44, 1).Value = "Afmelding nieuwsbrief"
$logs.Cells.Item(44, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(44, 3).Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Cells.Item(44, 4).Value = "Afmelding / Nieuwsbrief"
$logs.Cells.Item(44, 6).Value = "2025-06-19 19:50:27"
$logs.Cells.Item(44, 7).Value = "Nee"

# ---------------------------------------------------------------------------
# 2) Logs sheet: extend the conditional-formatting ranges to row 44
# ---------------------------------------------------------------------------
$dCond = $logs.Range("D2:D42").FormatConditions
for ($i = 1; $i -le $dCond.Count; $i++) {
    $dCond.Item($i).ModifyAppliesToRange($logs.Range("D2:D44"))
}

$gCond = $logs.Range("G2:G42").FormatConditions
for ($i = 1; $i -le $gCond.Count; $i++) {
    $gCond.Item($i).ModifyAppliesToRange($logs.Range("G2:G44"))
}

# ---------------------------------------------------------------------------
# 3) Dashboard sheet: update the category/count summary table
# ---------------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(6, 1).Value = "Offerte-aanvraag"
$dash.Cells.Item(6, 2).Value = 3

$dash.Cells.Item(7, 1).Value = "Bestelling"
$dash.Cells.Item(7, 2).Value = 3

$dash.Cells.Item(8, 1).Value = "Afmelding / Nieuwsbrief"
$dash.Cells.Item(8, 2).Value = 2

$dash.Cells.Item(9, 1).Value = "Openingstijden"
$dash.Cells.Item(9, 2).Value = 2

$dash.Cells.Item(10, 1).Value = "Informatieaanvraag"
$dash.Cells.Item(10, 2).Value = 1

$dash.Cells.Item(11, 1).Value = "Samenwerking"
$dash.Cells.Item(11, 2).Value = 1

$dash.Cells.Item(12, 1).Value = "Sollicitatie / Vacature"
$dash.Cells.Item(12, 2).Value = 1

$dash.Cells.Item(13, 1).Value = "Samenwerking / Partnerverzoek"
$dash.Cells.Item(13, 2).Value = 1

# ---------------------------------------------------------------------------
# 4) Dashboard chart: extend the category/value series ranges to row 13
# ---------------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$13,Dashboard!`$B`$2:`$B`$13,1)"
